# Append a new mapping row (row 16) that mirrors the
# "field_ddh_harvest_src" / "EnergyData.info" pair already used in row 2,
# but flagged FALSE the way row 15's repeated "field_wbddh_data_class"
# entry is -- i.e. an additional (machine_name, list_value_name, eex_value)
# resource-metadata mapping row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A16").Value = "field_ddh_harvest_src"
$ws.Range("B16").Value = "EnergyData.info"
$ws.Range("D16").Value = $false

# Touch the page setup (orientation) the way the source workbook's commit
# shows it was -- this is the only page-setup attribute that ends up
# persisted.
$ws.PageSetup.Orientation = 1

# Match the author's final selection captured in the sheet view.
$ws.Range("I18").Select()
